$d = $word.ActiveDocument

# Create the three new character styles
$sGaNStyle = $d.Styles.Add("GaNStyle", 2)
$sGaNStyle.Font.Name = "Calibri"
$sGaNStyle.Font.Size = 14

$sGaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$sGaNParagraph.Font.Name = "Calibri"
$sGaNParagraph.Font.Size = 10

$sGaNLinks = $d.Styles.Add("GaNLinks", 2)
$sGaNLinks.Font.Name = "Calibri"
$sGaNLinks.Font.Bold = $true
$sGaNLinks.Font.Color = 8388608
$sGaNLinks.Font.Size = 9.5
$sGaNLinks.Font.Underline = 1

# Apply GaNStyle to each of the 4 "2022 date" heading runs
$needle1 = "2022 Ημερομηνίες παρατήρησης για τον  Αστερισμός Ωρίωνα: 16-25 Ιανουαρίου, 14-23 Φεβρουαρίου, 14-24 Μαρτίου"
$r = $d.Content
$r.Start = 0
$r.End = $d.Content.End
while ($r.Find.Execute($needle1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $r.Style = "GaNStyle"
    $r.Collapse(0)
    $r.End = $d.Content.End
}

# Apply GaNParagraph to the campaign description paragraph run
$needle2 = "Συμμετέχετε σε μία παγκόσμια καμπάνια για να παρατηρήσετε και να καταγράψετε τη φωτεινότητα των πιο αμυδρά ορατών άστρων σαν μέσο για την μέτρηση της Φωτορρύπανσης σε μία δεδομένη περιοχή. Με τον εντοπισμό και την παρατήρηση του  Αστερισμός Ωρίωνα στον νυχτερινό ουρανό καθώς και με την σύγκριση των ανωτέρω με τα διαγράμματα για τα μεγέθη των άστρων,  άνθρωποι από όλον τον κόσμο θα μάθουν πώς τα φώτα στην κοινότητά τους συμβάλλουν στην Φωτορρύπανση. Με την κατάθεση των πορισμάτων τους στην ιστοσελίδα θα δημιουργηθεί ένα αρχείο σχετικά με το τι μπορεί να δει κανείς στον νυχτερινό ουρανό."
$r2 = $d.Content
$r2.Start = 0
$r2.End = $d.Content.End
if ($r2.Find.Execute($needle2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $r2.Style = "GaNParagraph"
}

# Apply GaNLinks to the Jan Hollan credit line run
$needle3 = "Τα διαγράμματα αυτού του αρχείου επιμελήθηκε ο Jan Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$r3 = $d.Content
$r3.Start = 0
$r3.End = $d.Content.End
if ($r3.Find.Execute($needle3, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $r3.Style = "GaNLinks"
}
